$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backward extension of the real-time data ---
# The sheet originally held 30 yearly rows (1995-2024) in rows 2-31.
# This change prepends 11 additional years (1984-1994) of history,
# pushing the existing rows down and extending the sheet through row 42.

# The new rows (32-42) created at the bottom of the sheet fall outside the
# old dimension (A1:E31), so they do not yet carry the date style used by
# the rest of column A (numFmt "YYYY-MM-DD HH:MM:SS" with borders/bold).
# Copy that formatting from an existing date cell before writing values.
$ws.Range("A2").Copy()
$ws.Range("A32:A42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Write out the full data table: the 11 newly-added historical rows
# followed by the pre-existing rows, now occupying rows 2 through 42.
$ws.Range("A2").Value = 31047
$ws.Range("B2").Value = 1984
$ws.Range("C2").Value = 2.833670241322217
$ws.Range("D2").Value = 1985
$ws.Range("E2").Value = 4.978977805976226
$ws.Range("A3").Value = 31412
$ws.Range("B3").Value = 1985
$ws.Range("C3").Value = 2.740628897120945
$ws.Range("D3").Value = 1986
$ws.Range("E3").Value = 4.840042388885646
$ws.Range("A4").Value = 31777
$ws.Range("B4").Value = 1986
$ws.Range("C4").Value = 2.269459987912947
$ws.Range("D4").Value = 1987
$ws.Range("E4").Value = 3.593781657196393
$ws.Range("A5").Value = 32142
$ws.Range("B5").Value = 1987
$ws.Range("C5").Value = 1.253514454810789
$ws.Range("D5").Value = 1988
$ws.Range("E5").Value = 5.084502077712005
$ws.Range("A6").Value = 32508
$ws.Range("B6").Value = 1988
$ws.Range("C6").Value = 3.509161092519553
$ws.Range("D6").Value = 1989
$ws.Range("E6").Value = 5.511076843601681
$ws.Range("A7").Value = 32873
$ws.Range("B7").Value = 1989
$ws.Range("C7").Value = 3.898460078540933
$ws.Range("D7").Value = 1990
$ws.Range("E7").Value = 2.951715842334024
$ws.Range("A8").Value = 33238
$ws.Range("B8").Value = 1990
$ws.Range("C8").Value = 5.356103277865332
$ws.Range("D8").Value = 1991
$ws.Range("E8").Value = 6.3181560832964
$ws.Range("A9").Value = 33603
$ws.Range("B9").Value = 1991
$ws.Range("C9").Value = 5.955905607167122
$ws.Range("D9").Value = 1992
$ws.Range("E9").Value = -0.005898890116151634
$ws.Range("A10").Value = 33969
$ws.Range("B10").Value = 1992
$ws.Range("C10").Value = 1.850401149566561
$ws.Range("D10").Value = 1993
$ws.Range("E10").Value = -0.4495646332120296
$ws.Range("A11").Value = 34334
$ws.Range("B11").Value = 1993
$ws.Range("C11").Value = -0.9857661435315745
$ws.Range("D11").Value = 1994
$ws.Range("E11").Value = 2.795029892345036
$ws.Range("A12").Value = 34699
$ws.Range("B12").Value = 1994
$ws.Range("C12").Value = 3.052254893522388
$ws.Range("D12").Value = 1995
$ws.Range("E12").Value = 3.383052772393214
$ws.Range("A13").Value = 35040
$ws.Range("B13").Value = 1995
$ws.Range("C13").Value = 2.234710814035812
$ws.Range("D13").Value = 1996
$ws.Range("E13").Value = 1.463127579670287
$ws.Range("A14").Value = 35403
$ws.Range("B14").Value = 1996
$ws.Range("C14").Value = 1.595002781738275
$ws.Range("D14").Value = 1997
$ws.Range("E14").Value = 4.207635715208324
$ws.Range("A15").Value = 35768
$ws.Range("B15").Value = 1997
$ws.Range("C15").Value = 2.499560583078497
$ws.Range("D15").Value = 1998
$ws.Range("E15").Value = 3.492506333467071
$ws.Range("A16").Value = 36132
$ws.Range("B16").Value = 1998
$ws.Range("C16").Value = 2.812603855740181
$ws.Range("D16").Value = 1999
$ws.Range("E16").Value = 2.334197296693863
$ws.Range("A17").Value = 36501
$ws.Range("B17").Value = 1999
$ws.Range("C17").Value = 1.188004848513446
$ws.Range("D17").Value = 2000
$ws.Range("E17").Value = 2.032004888754391
$ws.Range("A18").Value = 36858
$ws.Range("B18").Value = 2000
$ws.Range("C18").Value = 3.277038745546235
$ws.Range("D18").Value = 2001
$ws.Range("E18").Value = 3.09884301635126
$ws.Range("A19").Value = 37222
$ws.Range("B19").Value = 2001
$ws.Range("C19").Value = 0.7513248531724415
$ws.Range("D19").Value = 2002
$ws.Range("E19").Value = -0.408724114026926
$ws.Range("A20").Value = 37581
$ws.Range("B20").Value = 2002
$ws.Range("C20").Value = 0.2537741062064169
$ws.Range("D20").Value = 2003
$ws.Range("E20").Value = 0.9118162660485263
$ws.Range("A21").Value = 37938
$ws.Range("B21").Value = 2003
$ws.Range("C21").Value = -0.1535080579381121
$ws.Range("D21").Value = 2004
$ws.Range("E21").Value = 0.3435726964089891
$ws.Range("A22").Value = 38302
$ws.Range("B22").Value = 2004
$ws.Range("C22").Value = 1.171834509066594
$ws.Range("D22").Value = 2005
$ws.Range("E22").Value = 0.8394840956263971
$ws.Range("A23").Value = 38671
$ws.Range("B23").Value = 2005
$ws.Range("C23").Value = 1.120380359544382
$ws.Range("D23").Value = 2006
$ws.Range("E23").Value = 1.940699468213469
$ws.Range("A24").Value = 39035
$ws.Range("B24").Value = 2006
$ws.Range("C24").Value = 2.691354324129258
$ws.Range("D24").Value = 2007
$ws.Range("E24").Value = 3.187301687590338
$ws.Range("A25").Value = 39400
$ws.Range("B25").Value = 2007
$ws.Range("C25").Value = 2.652245539637632
$ws.Range("D25").Value = 2008
$ws.Range("E25").Value = 2.158031012958861
$ws.Range("A26").Value = 39765
$ws.Range("B26").Value = 2008
$ws.Range("C26").Value = 1.327195601304898
$ws.Range("D26").Value = 2009
$ws.Range("E26").Value = -1.941693908020603
$ws.Range("A27").Value = 40130
$ws.Range("B27").Value = 2009
$ws.Range("C27").Value = -4.803590807538871
$ws.Range("D27").Value = 2010
$ws.Range("E27").Value = 2.536922056245872
$ws.Range("A28").Value = 40494
$ws.Range("B28").Value = 2010
$ws.Range("C28").Value = 3.776429555840499
$ws.Range("D28").Value = 2011
$ws.Range("E28").Value = 5.124900822223233
$ws.Range("A29").Value = 40862
$ws.Range("B29").Value = 2011
$ws.Range("C29").Value = 3.167941427237042
$ws.Range("D29").Value = 2012
$ws.Range("E29").Value = 1.70423418303296
$ws.Range("A30").Value = 41228
$ws.Range("B30").Value = 2012
$ws.Range("C30").Value = 1.072335020576287
$ws.Range("D30").Value = 2013
$ws.Range("E30").Value = 0.990934028412549
$ws.Range("A31").Value = 41592
$ws.Range("B31").Value = 2013
$ws.Range("C31").Value = 0.5676944965793185
$ws.Range("D31").Value = 2014
$ws.Range("E31").Value = 1.859803271823757
$ws.Range("A32").Value = 41957
$ws.Range("B32").Value = 2014
$ws.Range("C32").Value = 1.417171832295883
$ws.Range("D32").Value = 2015
$ws.Range("E32").Value = 0.05331272828721367
$ws.Range("A33").Value = 42321
$ws.Range("B33").Value = 2015
$ws.Range("C33").Value = 1.475252114130599
$ws.Range("D33").Value = 2016
$ws.Range("E33").Value = 1.442973638880907
$ws.Range("A34").Value = 42689
$ws.Range("B34").Value = 2016
$ws.Range("C34").Value = 1.71887541289224
$ws.Range("D34").Value = 2017
$ws.Range("E34").Value = 1.076548192761484
$ws.Range("A35").Value = 43053
$ws.Range("B35").Value = 2017
$ws.Range("C35").Value = 2.581636142651922
$ws.Range("D35").Value = 2018
$ws.Range("E35").Value = 3.064375402422015
$ws.Range("A36").Value = 43418
$ws.Range("B36").Value = 2018
$ws.Range("C36").Value = 1.471137749280693
$ws.Range("D36").Value = 2019
$ws.Range("E36").Value = 0.1137080120319656
$ws.Range("A37").Value = 43783
$ws.Range("B37").Value = 2019
$ws.Range("C37").Value = 0.508332909595044
$ws.Range("D37").Value = 2020
$ws.Range("E37").Value = -0.1203207525434236
$ws.Range("A38").Value = 44159
$ws.Range("B38").Value = 2020
$ws.Range("C38").Value = -5.494775307949129
$ws.Range("D38").Value = 2021
$ws.Range("E38").Value = 8.235743591092737
$ws.Range("A39").Value = 44525
$ws.Range("B39").Value = 2021
$ws.Range("C39").Value = 3.149343082976164
$ws.Range("D39").Value = 2022
$ws.Range("E39").Value = 7.432336632701175
$ws.Range("A40").Value = 44890
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 1.995866057153428
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 1.189587957345273
$ws.Range("A41").Value = 45254
$ws.Range("B41").Value = 2023
$ws.Range("C41").Value = -0.1168430792840458
$ws.Range("D41").Value = 2024
$ws.Range("E41").Value = -0.1427298585871872
$ws.Range("A42").Value = 45618
$ws.Range("B42").Value = 2024
$ws.Range("C42").Value = -0.1775688094211469
$ws.Range("D42").Value = 2025
$ws.Range("E42").Value = -0.1265568156813002
